$wb = $excel.ActiveWorkbook

# OFF sheet (sheet1): row 2 ("H") - update Short/Deep Att/Comp/Int values
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 159
$wsOff.Range("C2").Value = 119
$wsOff.Range("D2").Value = 46
$wsOff.Range("E2").Value = 22
$wsOff.Range("F2").Value = 5

# DEF sheet (sheet2): row 2 ("H") - update Short/Deep Att/Comp/Int values
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 234
$wsDef.Range("C2").Value = 167
$wsDef.Range("D2").Value = 39
$wsDef.Range("E2").Value = 16
